{"js": "// Word JavaScript API (Office.js) edit script.\n// Body of: async (context) => { ... }\n//\n// Changes applied (per the commit \"add a techniques sheet too\"):\n//   1. \"Kicho E Jang\"   -> \"Kicho I Jang\"     (Forms table)\n//   2. \"Kicho Sam Jan\"  -> \"Kicho Sam Jang\"   (Forms table, missing trailing \"g\")\n//   3. \"Palgue E Jang\"  -> \"Palgue I Jang\"    (Forms table)\n//   4. \"Front\u2014Back: 1 Hand/1 Foot Technique\" -> \"Front-Back: 1 Hand/1 Foot Technique\"\n//      (em dash U+2014 replaced with a plain hyphen, Special/Authorization table)\n//\n// NOTE: the source text uses NON-BREAKING SPACES (U+00A0) between the words,\n// not regular spaces, so the literal search/replacement strings below use\n// \"\\u00a0\" to match exactly.\n\nconst NBSP = \"\\u00a0\";\n\nasync function replaceExact(body, findText, replaceText) {\n  const results = body.search(findText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${JSON.stringify(findText)}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(replaceText, \"Replace\");\n  }\n  await context.sync();\n}\n\nconst body = context.document.body;\n\nawait replaceExact(\n  body,\n  `Kicho${NBSP}E${NBSP}Jang`,\n  `Kicho${NBSP}I${NBSP}Jang`\n);\n\nawait replaceExact(\n  body,\n  `Kicho${NBSP}Sam${NBSP}Jan`,\n  `Kicho${NBSP}Sam${NBSP}Jang`\n);\n\nawait replaceExact(\n  body,\n  `Palgue${NBSP}E${NBSP}Jang`,\n  `Palgue${NBSP}I${NBSP}Jang`\n);\n\nawait replaceExact(\n  body,\n  `Front\\u2014Back:${NBSP}1${NBSP}Hand/1${NBSP}Foot${NBSP}Technique`,\n  `Front-Back:${NBSP}1${NBSP}Hand/1${NBSP}Foot${NBSP}Technique`\n);\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Changes applied (per the commit \"add a techniques sheet too\"):\n#   1. \"Kicho E Jang\"   -> \"Kicho I Jang\"     (Forms table)\n#   2. \"Kicho Sam Jan\"  -> \"Kicho Sam Jang\"   (Forms table, missing trailing \"g\")\n#   3. \"Palgue E Jang\"  -> \"Palgue I Jang\"    (Forms table)\n#   4. \"Front\u2014Back: 1 Hand/1 Foot Technique\" -> \"Front-Back: 1 Hand/1 Foot Technique\"\n#      (em dash U+2014 replaced with a plain hyphen, Special/Authorization table)\n#\n# NOTE: the source text uses NON-BREAKING SPACES (U+00A0) between the words,\n# not regular spaces, so the literal search/replacement strings below use\n# [char]0x00A0 to match exactly.\n\n$d = $word.ActiveDocument\n\n$NBSP = [char]0x00A0\n$EMDASH = [char]0x2014\n\nfunction Replace-ExactText($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $found = $find.Execute($findText, $false, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n}\n\nReplace-ExactText \"Kicho${NBSP}E${NBSP}Jang\" \"Kicho${NBSP}I${NBSP}Jang\"\nReplace-ExactText \"Kicho${NBSP}Sam${NBSP}Jan\" \"Kicho${NBSP}Sam${NBSP}Jang\"\nReplace-ExactText \"Palgue${NBSP}E${NBSP}Jang\" \"Palgue${NBSP}I${NBSP}Jang\"\nReplace-ExactText \"Front${EMDASH}Back:${NBSP}1${NBSP}Hand/1${NBSP}Foot${NBSP}Technique\" \"Front-Back:${NBSP}1${NBSP}Hand/1${NBSP}Foot${NBSP}Technique\"\n"}
